# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.686.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.598.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.598.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.663.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0753"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0516"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.289.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.26%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  +17.26%  "
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.735.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
